$d = $word.ActiveDocument

# Locate the paragraph that holds the hidden "_GoBack" bookmark (it is an
# otherwise content-less paragraph near the end of the document). The new
# "Cycle 4" user stories are inserted immediately before it, and one extra
# blank paragraph is inserted immediately after it.

function Find-BookmarkParagraphIndex($doc) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.WordOpenXML -like "*_GoBack*") {
            return $i
        }
    }
    return -1
}

$bookmarkIndex = Find-BookmarkParagraphIndex $d
$bookmarkPara = $d.Paragraphs.Item($bookmarkIndex)

$cycle4Text = "`r" + `
    "Cycle 4:`r" + `
    "As a player, I would like to be able to play with an Xbox 360 controller.`r" + `
    "As a player, I would like to see a projectile.`r" + `
    "As a player, I would like to have the heroes be able to switch places.`r" + `
    "As a player I would like to have a load screen where the best past scores are displayed with a profile name.`r" + `
    "As a player, I would like to play a survival mode in order to be further challenged.`r" + `
    "As a player, I would like random events such as new enemies spawning in random areas or new emerald and gold bags being added to the level.`r" + `
    "As a player, I would like to the screen to have a different image for moving in each direction in order to clearly see which direction the player is facing.`r" + `
    "As a player, I would like to hear sounds and music while playing.`r"

$insertPoint = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)
$insertPoint.InsertBefore($cycle4Text)

# Add one blank paragraph right after the bookmark paragraph, before the
# document's pre-existing trailing blank paragraph (which is always last).
$trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint2 = $d.Range($trailingPara.Range.Start, $trailingPara.Range.Start)
$insertPoint2.InsertBefore("`r")
